# Adding CSV support: new cascading "region" -> "country" select (backed by
# regions.csv / countries_csv queries) plus a hint for the existing content
# provider example, and a hint/choice_filter column pair on the survey sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# survey sheet
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# New columns: H = hint, I = choice_filter
$survey.Range("H1").Value = "hint"
$survey.Range("I1").Value = "choice_filter"

# Make room for the two new rows (region / country) right above the
# "content_provider_test" row.
$survey.Rows.Item(17).Insert()
$survey.Rows.Item(17).Insert()

$survey.Range("B17").Value = "select_one regions_csv"
$survey.Range("F17").Value = "region"
$survey.Range("G17").Value = "Choose a region:"

$survey.Range("B18").Value = "select_one countries_csv"
$survey.Range("C18").Value = "dropdown"
$survey.Range("F18").Value = "country"
$survey.Range("G18").Value = "Choose a country:"
$survey.Range("I18").Value = "_.where(context, {`n  region: data('region')`n})"

# cp_test row (now at row 19) gets a new hint
$survey.Range("H19").Value = "You will need to install a content provider app for the query to work. There is an example app available here: https://github.com/nathanathan/FileContentProviderExample"

# ---------------------------------------------------------------------------
# queries sheet
# ---------------------------------------------------------------------------
$queries = $wb.Worksheets.Item("queries")

# Reuse the blank row 4 / insert one more so both new queries fit above the
# (updated) content_provider_test query.
$queries.Rows.Item(4).Insert()

$queries.Range("A4").Value = "countries_csv"
$queries.Range("B4").Value = """regions.csv"""
$queries.Range("C4").Value = "_.chain(context).pluck('region').uniq().map(function(region){`nreturn {name:region, label:region};`n})"

$queries.Range("A5").Value = "regions_csv"
$queries.Range("B5").Value = """regions.csv"""
$queries.Range("C5").Value = "_.map(context, function(place){`nplace.name = place.country;`nplace.label = place.country;`nreturn place;`n})"

$queries.Range("A6").Value = "content_provider_test"
$queries.Range("B6").Value = """content://org.opendatakit.FileContentProviderExample/"""
$queries.Range("C6").Value = "[context]"

# Widen the callback column to fit the new (longer) queries.
$queries.Columns.Item(3).ColumnWidth = 42.8
